$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold prices formatted as plain text (e.g. "307.43").
# Several of the updated prices look like ordinary numbers, so without
# forcing a text number format first, Excel would silently reinterpret
# them as numeric values (losing trailing zeros / exact formatting).
$textCells = @(
    "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D23", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the refreshed coin data (price + 1h volume change, and the
# FraxShare / TheSandbox row swap).
$ws.Range("D2").Value = "26.327.51"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.792.82"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "307.43"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("D7").Value = "0.4528"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("D8").Value = "0.3592"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "45.45"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "0.07091"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "0.8834"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "0.07825"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "19.46"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "1.835.10"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "5.287"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "6.330"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "84.65"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "1.009"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "0.000008532"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").Value = "26.358.60"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").Value = "4.985"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "2.030.27"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").Value = "1.976"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "152.19"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").Value = "17.86"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").Value = "2.028"
$ws.Range("E29").Value = "  +3.42%  "
$ws.Range("D30").Value = "111.91"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "0.08678"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").Value = "3.063"
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").Value = "4.442"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "0.7246"
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("D36").Value = "2.714"
$ws.Range("E36").Value = "  +6.02%  "
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("D38").Value = "1.071"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").Value = "0.01928"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "0.05103"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "2.877"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "6.860"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.5052"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("D44").Value = "0.1512"
$ws.Range("E44").Value = "  -5.19%  "
$ws.Range("D45").Value = "7.985"
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("D46").Value = "1.008"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "0.4627"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").Value = "100.89"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").Value = "9.812"
$ws.Range("E49").Value = "  -3.57%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("E51").Value = "  -2.08%  "

# Restore the default (unstyled) cell format now that the text values
# are safely stored, so formatting matches the original workbook.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).ClearFormats()
}
